$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 615.74194
$ws.Range("I28").Value = 369.21738
$ws.Range("J28").Value = 1324.5
$ws.Range("K28").Value = 369.21738
$ws.Range("L28").Value = 1324.5
$ws.Range("M28").Value = 115.78262
$ws.Range("N28").Value = -2294.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1753.5
$ws.Range("I26").Value = 1753.5
$ws.Range("K26").Value = 1753.5
$ws.Range("M26").Value = -1423.5
$ws.Range("H122").Value = 20837868
$ws.Range("I122").Value = 5956
$ws.Range("J122").Value = 31253824
$ws.Range("K122").Value = 17868
$ws.Range("L122").Value = 93761472
$ws.Range("M122").Value = -15418
$ws.Range("N122").Value = -93766372

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 1552.1111
$ws.Range("I25").Value = 1308.625
$ws.Range("K25").Value = 1308.625
$ws.Range("M25").Value = -1073.625
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H86").Value = 2118.976
$ws.Range("J86").Value = 1751.75
$ws.Range("L86").Value = 1751.75
$ws.Range("N86").Value = -3997.75
$ws.Range("H89").Value = 2118.976
$ws.Range("J89").Value = 1751.75
$ws.Range("L89").Value = 8758.75
$ws.Range("N89").Value = -19990.75
$ws.Range("H94").Value = 1928.6316
$ws.Range("I94").Value = 1566.8182
$ws.Range("J94").Value = 2426.125
$ws.Range("K94").Value = 1566.8182
$ws.Range("L94").Value = 2426.125
$ws.Range("M94").Value = -1115.8182
$ws.Range("N94").Value = -3328.125
$ws.Range("H132").Value = 60060.375
$ws.Range("J132").Value = 60060.375
$ws.Range("L132").Value = 60060.375
$ws.Range("N132").Value = -70180.375
$ws.Range("H134").Value = 6341
$ws.Range("I134").Value = 6553.1113
$ws.Range("K134").Value = 19659.3339
$ws.Range("M134").Value = -17124.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 32000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 32000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 32000
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -33120

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 6000
$ws.Range("J43").Value = 6000
$ws.Range("L43").Value = 18000
$ws.Range("N43").Value = -18228
$ws.Range("H46").Value = 2812.125
$ws.Range("J46").Value = 2979.6
$ws.Range("L46").Value = 8938.799999999999
$ws.Range("N46").Value = -9120.799999999999
$ws.Range("H55").Value = 4985
$ws.Range("J55").Value = 4985
$ws.Range("L55").Value = 14955
$ws.Range("N55").Value = -15309
$ws.Range("H76").Value = 3000
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 9000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -8617
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3000
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -7674
$ws.Range("N79").ClearContents()
$ws.Range("H94").Value = 3800
$ws.Range("J94").Value = 3800
$ws.Range("L94").Value = 11400
$ws.Range("N94").Value = -12752
$ws.Range("H109").Value = 2590
$ws.Range("I109").Value = 700
$ws.Range("J109").Value = 2716
$ws.Range("K109").Value = 2100
$ws.Range("L109").Value = 8148
$ws.Range("M109").Value = -1060
$ws.Range("N109").Value = -10228
$ws.Range("H112").Value = 4875.8237
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 4875.8237
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 14627.4711
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -16843.4711
$ws.Range("H125").Value = 2746.353
$ws.Range("I125").Value = 844
$ws.Range("K125").Value = 2532
$ws.Range("M125").Value = 2388
$ws.Range("H131").Value = 60355.25
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 60355.25
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 181065.75
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -191145.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 14400
$ws.Range("J15").Value = 14400
$ws.Range("L15").Value = 14400
$ws.Range("N15").Value = -14976
$ws.Range("H81").Value = 14400
$ws.Range("J81").Value = 14400
$ws.Range("L81").Value = 14400
$ws.Range("N81").Value = -16396
$ws.Range("H84").Value = 14400
$ws.Range("J84").Value = 14400
$ws.Range("L84").Value = 43200
$ws.Range("N84").Value = -53184

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1041.1666
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 1049.4
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 1049.4
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -1425.4
$ws.Range("H82").Value = 1576.2916
$ws.Range("I82").Value = 1205.5385
$ws.Range("J82").Value = 2014.4546
$ws.Range("K82").Value = 1205.5385
$ws.Range("L82").Value = 2014.4546
$ws.Range("M82").Value = -844.5385000000001
$ws.Range("N82").Value = -2736.4546
$ws.Range("H85").Value = 1576.2916
$ws.Range("I85").Value = 1205.5385
$ws.Range("J85").Value = 2014.4546
$ws.Range("K85").Value = 1205.5385
$ws.Range("L85").Value = 2014.4546
$ws.Range("M85").Value = 42.46149999999989
$ws.Range("N85").Value = -4510.4546
$ws.Range("H93").Value = 2617.6667
$ws.Range("I93").Value = 2617.6667
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2617.6667
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1369.6667
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3226.7368
$ws.Range("I62").Value = 3698
$ws.Range("J62").Value = 3058.4285
$ws.Range("K62").Value = 3698
$ws.Range("L62").Value = 3058.4285
$ws.Range("M62").Value = -3074
$ws.Range("N62").Value = -4306.4285
$ws.Range("H65").Value = 3226.7368
$ws.Range("I65").Value = 3698
$ws.Range("J65").Value = 3058.4285
$ws.Range("K65").Value = 18490
$ws.Range("L65").Value = 15292.1425
$ws.Range("M65").Value = -15370
$ws.Range("N65").Value = -21532.1425
$ws.Range("H81").Value = 2680
$ws.Range("I81").Value = 2133.3333
$ws.Range("J81").Value = 3500
$ws.Range("K81").Value = 4266.6666
$ws.Range("L81").Value = 7000
$ws.Range("M81").Value = -3205.6666
$ws.Range("N81").Value = -9122
$ws.Range("H84").Value = 2680
$ws.Range("I84").Value = 2133.3333
$ws.Range("J84").Value = 3500
$ws.Range("K84").Value = 21333.333
$ws.Range("L84").Value = 35000
$ws.Range("M84").Value = -16029.333
$ws.Range("N84").Value = -45608
